$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Datos actualizados..." timestamp note (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 00:05"

# --- 2. Update numeric stats for countries whose rank/position didn't change ---

# Estados Unidos
$r = $ws.Columns("A").Find("Estados Unidos").Row
$ws.Cells.Item($r, 2).Value = 1742443
$ws.Cells.Item($r, 3).Value = 17168
$ws.Cells.Item($r, 4).Value = 485929
$ws.Cells.Item($r, 5).Value = 1154667
$ws.Cells.Item($r, 7).Value = 1275
$ws.Cells.Item($r, 8).Value = 101847

# Peru
$r = $ws.Columns("A").Find("Peru").Row
$ws.Cells.Item($r, 2).Value = 135905
$ws.Cells.Item($r, 3).Value = 6154
$ws.Cells.Item($r, 4).Value = 52906
$ws.Cells.Item($r, 5).Value = 79016
$ws.Cells.Item($r, 7).Value = 195
$ws.Cells.Item($r, 8).Value = 3983

# Venezuela
$r = $ws.Columns("A").Find("Venezuela").Row
$ws.Cells.Item($r, 2).Value = 1245
$ws.Cells.Item($r, 3).Value = 34
$ws.Cells.Item($r, 5).Value = 932

# --- 3. Guayana Francesa moved up in the ranking, overtaking Togo and Cabo
#        Verde (whose own figures are unchanged, they just shift down a row) ---

# Remove the old "Guayana Francesa" row (Togo/Cabo Verde/... shift up by one)
$gfRow = $ws.Columns("A").Find("Guayana Francesa").Row
$ws.Rows($gfRow).Delete()

# Recompute Togo's row (it moved up by one after the delete) and insert a new
# row right before it for the updated "Guayana Francesa" entry
$togoRow = $ws.Columns("A").Find("Togo").Row
$ws.Rows($togoRow).Insert()

$ws.Cells.Item($togoRow, 1).Value = "Guayana Francesa"
$ws.Cells.Item($togoRow, 2).Value = 406
$ws.Cells.Item($togoRow, 3).Value = 22
$ws.Cells.Item($togoRow, 4).Value = 150
$ws.Cells.Item($togoRow, 5).Value = 255
$ws.Cells.Item($togoRow, 6).Value = 0
$ws.Cells.Item($togoRow, 7).Value = 0
$ws.Cells.Item($togoRow, 8).Value = 1
